# Fix the "fixed asset schedule": cells that previously held empty string
# placeholders in the Salvage Value column (N) and, for the later
# straight-line/blanket depreciation rows, also in the Years used(sold
# items) / Accumulated Depreciation (on Sold items) / Net columns (Q, W, X, Y)
# should instead hold an explicit numeric 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Salvage Value column (N) for every data row, 2 through 334
$ws.Range("N2:N334").Value = 0

# Years used(sold items) / Accumulated Depreciation on Sold items /
# Accumulated Depreciation Net / Book Value columns for rows 302-334
$ws.Range("Q302:Q334").Value = 0
$ws.Range("W302:W334").Value = 0
$ws.Range("X302:X334").Value = 0
$ws.Range("Y302:Y334").Value = 0
